$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 98, shifting existing rows 98-129 down to 99-130.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new record's data.
$ws.Cells.Item(98, 1).Value = 5
$ws.Cells.Item(98, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(98, 3).Value = "Maule"
$ws.Cells.Item(98, 4).Value = 45215
$ws.Cells.Item(98, 5).Value = 7
$ws.Cells.Item(98, 6).Value = 100112026
$ws.Cells.Item(98, 7).Value = "Haba"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 700
$ws.Cells.Item(98, 11).Value = 10000
$ws.Cells.Item(98, 12).Value = 11000
$ws.Cells.Item(98, 13).Value = 10286
$ws.Cells.Item(98, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(98, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(98, 16).Value = 411
$ws.Cells.Item(98, 17).Value = 25
$ws.Cells.Item(98, 18).Value = "Hortaliza"
